# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the old N/O/P ("Late" / heading / "Outstanding") columns
#   right to O/P/Q. The new column inherits the width of the column to its
#   left (column M) and is left blank (a spacer column), matching the
#   published workbook.
# - The "Repayment schedule" sheet becomes the active sheet/tab (it was
#   "NewLoanInput" before), with cell T6 selected on it.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N; existing N/O/P shift to O/P/Q, carrying
# their values/styles with them.
$ws3.Columns("N").Insert()

# Match the new column's width to its left-hand neighbour (column M), as
# Excel does visually when a column is inserted.
$ws3.Columns("N").ColumnWidth = $ws3.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab with T6 selected.
$ws3.Activate() | Out-Null
$ws3.Range("T6").Select() | Out-Null
